# "Criteria of success" text box on slide 1 (Google Shape;35;p1).
# The paragraph's trailing word "based " is removed from the end of the
# last run. We rebuild the run boundaries using Characters() sub-ranges
# so the resulting runs/formatting match the authored edit:
#   run2 absorbs run3's text plus run4's text up to "(future facilities"
#   run3 becomes ") " (keeping run3's original formatting)
#   run4 (now empty) collapses away.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$shp = $s.Shapes.Item("Google Shape;35;p1")
$tr  = $shp.TextFrame.TextRange

# Run 2 (" Present the list of facilities visitors would be willing to pay ")
# grows to include the old run3 text and most of run4's text.
$r2 = $tr.Characters(103, 65)
$r2.Text = " Present the list of facilities visitors would be willing to pay more. Capitalize from the resort’s existing facilities and present new investment plan (future facilities"

# Old run3 ("more. Capitalize ") becomes ") ".
$r3 = $tr.Characters(273, 17)
$r3.Text = ") "

# Old run4 ("from the resort’s existing facilities and present new
# investment plan (future facilities) based ") is now fully absorbed /
# superseded, so it collapses to nothing.
$r4 = $tr.Characters(275, 96)
$r4.Text = ""
